$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9972865504815053
$ws.Range("E2").Value = 0.9972865504815053

# Row 3
$ws.Range("D3").Value = 0.160986278550325
$ws.Range("E3").Value = 0.160986278550325

# Row 4
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.9997263831940452
$ws.Range("E4").Value = 0.9997263831940452

# Row 5
$ws.Range("D5").Value = [double]"7.949644102596584E-15"
$ws.Range("E5").Value = [double]"7.949644102596584E-15"

# Row 6
$ws.Range("D6").Value = 0.08582694435864666
$ws.Range("E6").Value = 0.08582694435864666

# Row 7
$ws.Range("D7").Value = 0.2305722876886664
$ws.Range("E7").Value = 0.7694277123113337

# Row 8
$ws.Range("D8").Value = 0.1140232415518011
$ws.Range("E8").Value = 0.8859767584481989

# Row 10
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0.6259873946277311
$ws.Range("E10").Value = 0.3740126053722689

# Row 11
$ws.Range("D11").Value = 0.9766296727181878
$ws.Range("E11").Value = 0.02337032728181221
$ws.Range("F11").Value = 1.850919723510742
